$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new motor entry (HobbyKing Turnigy Aerodrive SK3 5045, 450kv,
# stall speed unknown -> "?") as a new last row of the data block.
$ws.Rows(12).Insert()
$ws.Cells.Item(12,5).Value = "https://hobbyking.com/en_us/turnigy-aerodrive-sk3-5045-450kv-brushless-outrunner-motor.html?___store=en_us#qa[bW9kZT03JnBhZ2U9MSZxdWVzdGlvbl9zZWFyY2hfY29udGVudD0mcT0zMTk0]"
$ws.Cells.Item(12,1).Value = "?"
$ws.Cells.Item(12,2).Value = 1260
$ws.Cells.Item(12,3).Value = 275
$ws.Cells.Item(12,4).Value = 450

# Re-sort the whole data block ascending by Battery Power (W), column B,
# now that the new (heavier motor sizing) row is included.
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("B2:B12"))
$ws.Sort.SetRange($ws.Range("A2:E12"))
$ws.Sort.Apply()

# Match the post-edit selection/view state.
$ws.Range("E15").Select()
